# Updated cryptos list on Wed Oct  2 05:39:38 UTC 2024 with GitHub Actions
#
# Refreshes Price / Volume(1h) figures for the crypto ranking sheet and
# restores the correct row order for Toncoin / TRON (rows 11 and 12 had
# been swapped). All affected cells hold plain text (not numeric) values
# in the source workbook, so for cells whose new text looks like a number
# we briefly force a text NumberFormat while assigning the value and then
# clear the formatting again, which keeps the cell's value as text without
# leaving any extra formatting/style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.444.30"
$ws.Range("E2").Value = "  -3.59%  "

$ws.Range("D3").Value = "2.480.92"
$ws.Range("E3").Value = "  -6.01%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.39"
$ws.Range("D5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.48"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.68%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -3.30%  "

$ws.Range("D9").Value = "2.479.96"
$ws.Range("E9").Value = "  -5.88%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.42%  "

$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.44"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -6.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.357"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -6.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.18"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -7.79%  "

$ws.Range("D15").Value = "2.928.61"
$ws.Range("E15").Value = "  -5.97%  "

$ws.Range("E16").Value = "  -8.80%  "

$ws.Range("D17").Value = "61.346.13"
$ws.Range("E17").Value = "  -3.67%  "

$ws.Range("D18").Value = "2.488.21"
$ws.Range("E18").Value = "  -5.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.18"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -7.78%  "

$ws.Range("E20").Value = "  -7.71%  "

$ws.Range("E21").Value = "  -7.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "321.98"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -6.52%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("E24").Value = "  +0.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.12"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -5.73%  "

$ws.Range("D26").Value = "0.0₃0996"
$ws.Range("E26").Value = "  -8.57%  "

$ws.Range("D27").Value = "2.609.50"
$ws.Range("E27").Value = "  -5.71%  "

$ws.Range("E28").Value = "  -5.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "545.12"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -9.69%  "

$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.34"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -9.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.79"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.13%  "

$ws.Range("E33").Value = "  -5.52%  "

$ws.Range("E34").Value = "  -7.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.89"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -10.14%  "

$ws.Range("E37").Value = "  -10.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.383"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.95%  "

$ws.Range("E40").Value = "  -5.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "147.28"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.74"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -7.98%  "

$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.40"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.57%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.37"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -7.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "147.60"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -8.39%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.63"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.11"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -12.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0540"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -7.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.597"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.55%  "

$ws.Range("E51").Value = "  -4.78%  "
